# Fruta / hortaliza, semanal
# Insert a new daily record row at row 948 (pushing all existing rows
# from 948 downward by one), then populate the new row with the
# latest day's data for "Vega Modelo de Temuco" - Mandarina - Murcott.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 948; everything from the old
# row 948 onward (through the old row 1040) shifts down to 949..1041.
$ws.Rows.Item(948).Insert()

# Populate the newly inserted row 948 with the new observation.
$ws.Cells.Item(948, 1).Value  = 10
$ws.Cells.Item(948, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(948, 3).Value  = "La Araucanía"
$ws.Cells.Item(948, 4).Value  = 45166
$ws.Cells.Item(948, 5).Value  = 9
$ws.Cells.Item(948, 6).Value  = "Fruta"
$ws.Cells.Item(948, 7).Value  = 100102
$ws.Cells.Item(948, 8).Value  = "Cítricos"
$ws.Cells.Item(948, 9).Value  = 100102004
$ws.Cells.Item(948, 10).Value = "Mandarina"
$ws.Cells.Item(948, 11).Value = "Murcott"
$ws.Cells.Item(948, 12).Value = "Primera"
$ws.Cells.Item(948, 13).Value = 280
$ws.Cells.Item(948, 14).Value = 23000
$ws.Cells.Item(948, 15).Value = 23000
$ws.Cells.Item(948, 16).Value = 23000
$ws.Cells.Item(948, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(948, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(948, 19).Value = 1150
$ws.Cells.Item(948, 20).Value = 20
